$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 733.3889
$ws.Range("I28").Value = 285.16666
$ws.Range("J28").Value = 1629.8334
$ws.Range("K28").Value = 285.16666
$ws.Range("L28").Value = 1629.8334
$ws.Range("M28").Value = 199.83334
$ws.Range("N28").Value = -2599.8334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2949.1667
$ws.Range("I69").Value = 700
$ws.Range("J69").Value = 3399
$ws.Range("K69").Value = 2100
$ws.Range("L69").Value = 10197
$ws.Range("M69").Value = -1226
$ws.Range("N69").Value = -11945

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 2949.1667
$ws.Range("I72").Value = 700
$ws.Range("J72").Value = 3399
$ws.Range("K72").Value = 6300
$ws.Range("L72").Value = 30591
$ws.Range("M72").Value = -1932
$ws.Range("N72").Value = -39327

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 313.26315
$ws.Range("I92").Value = 293.18182
$ws.Range("J92").Value = 340.875
$ws.Range("K92").Value = 293.18182
$ws.Range("L92").Value = 340.875
$ws.Range("M92").Value = 954.81818
$ws.Range("N92").Value = -2836.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1065.4445
$ws.Range("J125").Value = 2349.5
$ws.Range("L125").Value = 21145.5
$ws.Range("N125").Value = -26065.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1500
$ws.Range("I127").Value = 1000
$ws.Range("J127").Value = 2000
$ws.Range("K127").Value = 3000
$ws.Range("L127").Value = 6000
$ws.Range("M127").Value = 1960
$ws.Range("N127").Value = -15920

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 373976.94
$ws.Range("I132").Value = 458667.2
$ws.Range("J132").Value = 1339.8
$ws.Range("K132").Value = 1376001.6
$ws.Range("L132").Value = 4019.4
$ws.Range("M132").Value = -1373471.6
$ws.Range("N132").Value = -9079.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 12000
$ws.Range("I38").Value = 6000
$ws.Range("J38").Value = 18000
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 18000
$ws.Range("M38").Value = -5533
$ws.Range("N38").Value = -18934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 250874.88
$ws.Range("I45").Value = 286399.84
$ws.Range("K45").Value = 286399.84
$ws.Range("M45").Value = -286022.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1939.8379
$ws.Range("I132").Value = 1564.174
$ws.Range("J132").Value = 2557
$ws.Range("K132").Value = 4692.522
$ws.Range("L132").Value = 7671
$ws.Range("M132").Value = -2162.522
$ws.Range("N132").Value = -12731

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 51425.25
$ws.Range("I134").Value = 67928.22
$ws.Range("J134").Value = 1916.3334
$ws.Range("K134").Value = 203784.66
$ws.Range("L134").Value = 5749.0002
$ws.Range("M134").Value = -201249.66
$ws.Range("N134").Value = -10819.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 908.2222
$ws.Range("I35").Value = 908.2222
$ws.Range("K35").Value = 908.2222
$ws.Range("M35").Value = -614.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 830.17645
$ws.Range("I105").Value = 817.2727
$ws.Range("J105").Value = 853.8333
$ws.Range("K105").Value = 817.2727
$ws.Range("L105").Value = 853.8333
$ws.Range("M105").Value = 929.7273
$ws.Range("N105").Value = -4347.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 14900
$ws.Range("J112").Value = 14900
$ws.Range("L112").Value = 14900
$ws.Range("N112").Value = -17854

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 23413.889
$ws.Range("J133").Value = 23413.889
$ws.Range("L133").Value = 23413.889
$ws.Range("N133").Value = -28473.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 35354950
$ws.Range("J96").Value = 35354950
$ws.Range("L96").Value = 106064850
$ws.Range("N96").Value = -106068968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1032.5
$ws.Range("J98").Value = 1169.5714
$ws.Range("L98").Value = 3508.7142
$ws.Range("N98").Value = -6504.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 7848.5
$ws.Range("J101").Value = 7848.5
$ws.Range("L101").Value = 23545.5
$ws.Range("N101").Value = -28413.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3055.5557
$ws.Range("J110").Value = 3616.6667
$ws.Range("L110").Value = 10850.0001
$ws.Range("N110").Value = -19030.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 483.63333
$ws.Range("I113").Value = 468.4375
$ws.Range("J113").Value = 501
$ws.Range("K113").Value = 1405.3125
$ws.Range("L113").Value = 1503
$ws.Range("M113").Value = 764.6875
$ws.Range("N113").Value = -5843

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 637
$ws.Range("I129").Value = 463.33334
$ws.Range("J129").Value = 785.8570999999999
$ws.Range("K129").Value = 1390.00002
$ws.Range("L129").Value = 2357.5713
$ws.Range("M129").Value = 3609.99998
$ws.Range("N129").Value = -12357.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1494975.6
$ws.Range("I131").Value = 4179.4614
$ws.Range("J131").Value = 1853871
$ws.Range("K131").Value = 12538.3842
$ws.Range("L131").Value = 5561613
$ws.Range("M131").Value = -7498.3842
$ws.Range("N131").Value = -5571693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 37982.17
$ws.Range("I137").Value = 2004.5454
$ws.Range("J137").Value = 59968.5
$ws.Range("K137").Value = 6013.6362
$ws.Range("L137").Value = 179905.5
$ws.Range("M137").Value = -913.6361999999999
$ws.Range("N137").Value = -190105.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 6747.8
$ws.Range("J36").Value = 4568.3335
$ws.Range("L36").Value = 4568.3335
$ws.Range("N36").Value = -5538.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1202
$ws.Range("J102").Value = 1000
$ws.Range("L102").Value = 1000
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2288.2683
$ws.Range("I132").Value = 1815.75
$ws.Range("J132").Value = 2955.353
$ws.Range("K132").Value = 5447.25
$ws.Range("L132").Value = 8866.059000000001
$ws.Range("M132").Value = -2917.25
$ws.Range("N132").Value = -13926.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1419
$ws.Range("I31").Value = 838
$ws.Range("K31").Value = 838
$ws.Range("M31").Value = -590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1981.2858
$ws.Range("I40").Value = 1931.4375
$ws.Range("J40").Value = 2140.8
$ws.Range("K40").Value = 1931.4375
$ws.Range("L40").Value = 2140.8
$ws.Range("M40").Value = -1795.4375
$ws.Range("N40").Value = -2412.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 51000
$ws.Range("I136").Value = 51000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 153000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -150450
$ws.Range("N136").ClearContents()

Write-Output "Applied all updates"